$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6881063333333334
$ws.Range("H2").Value = 2.064319
$ws.Range("I2").Value = 0.04678220357266529
$ws.Range("J2").Value = 0.04678220357266529
$ws.Range("M2").Value = 25.37147633333333
$ws.Range("N2").Value = 76.114429
$ws.Range("O2").Value = 0.5780881462719274
$ws.Range("P2").Value = 0.5780881462719274
$ws.Range("Q2").Value = 17.45827355098345
$ws.Range("R2").Value = 157.124461958851
$ws.Range("S2").Value = 0.02704423734183802
$ws.Range("T2").Value = 0.02704423734183802
$ws.Range("G3").Value = 0.6881063333333334
$ws.Range("H3").Value = 2.064319
$ws.Range("I3").Value = 0.04678220357266529
$ws.Range("J3").Value = 0.04678220357266529
$ws.Range("O3").Value = 0.2328552951919536
$ws.Range("P3").Value = 0.2328552951919536
$ws.Range("Q3").Value = 7.032234560547224
$ws.Range("R3").Value = 63.29011104492501
$ws.Range("S3").Value = 0.01089348382264304
$ws.Range("T3").Value = 0.01089348382264304
$ws.Range("G4").Value = 0.6881063333333334
$ws.Range("H4").Value = 2.064319
$ws.Range("I4").Value = 0.04678220357266529
$ws.Range("J4").Value = 0.04678220357266529
$ws.Range("O4").Value = 0.189056558536119
$ws.Range("P4").Value = 0.189056558536119
$ws.Range("Q4").Value = 5.709511839702222
$ws.Range("R4").Value = 51.38560655732
$ws.Range("S4").Value = 0.00884448240818423
$ws.Range("T4").Value = 0.00884448240818423
$ws.Range("I5").Value = 0.5086503334983149
$ws.Range("J5").Value = 0.5086503334983149
$ws.Range("M5").Value = 25.37147633333333
$ws.Range("N5").Value = 76.114429
$ws.Range("O5").Value = 0.5780881462719274
$ws.Range("P5").Value = 0.5780881462719274
$ws.Range("Q5").Value = 189.8191189352438
$ws.Range("R5").Value = 1708.372070417194
$ws.Range("S5").Value = 0.2940447283926385
$ws.Range("T5").Value = 0.2940447283926385
$ws.Range("I6").Value = 0.5086503334983149
$ws.Range("J6").Value = 0.5086503334983149
$ws.Range("O6").Value = 0.2328552951919536
$ws.Range("P6").Value = 0.2328552951919536
$ws.Range("S6").Value = 0.1184419235562358
$ws.Range("T6").Value = 0.1184419235562358
$ws.Range("I7").Value = 0.5086503334983149
$ws.Range("J7").Value = 0.5086503334983149
$ws.Range("O7").Value = 0.189056558536119
$ws.Range("P7").Value = 0.189056558536119
$ws.Range("S7").Value = 0.0961636815494406
$ws.Range("T7").Value = 0.09616368154944062
$ws.Range("I8").Value = 0.4445674629290199
$ws.Range("J8").Value = 0.4445674629290199
$ws.Range("M8").Value = 25.37147633333333
$ws.Range("N8").Value = 76.114429
$ws.Range("O8").Value = 0.5780881462719274
$ws.Range("P8").Value = 0.5780881462719274
$ws.Range("Q8").Value = 165.9045488873995
$ws.Range("R8").Value = 1493.140939986595
$ws.Range("S8").Value = 0.256999180537451
$ws.Range("T8").Value = 0.2569991805374509
$ws.Range("I9").Value = 0.4445674629290199
$ws.Range("J9").Value = 0.4445674629290199
$ws.Range("O9").Value = 0.2328552951919536
$ws.Range("P9").Value = 0.2328552951919536
$ws.Range("S9").Value = 0.1035198878130748
$ws.Range("T9").Value = 0.1035198878130748
$ws.Range("I10").Value = 0.4445674629290199
$ws.Range("J10").Value = 0.4445674629290199
$ws.Range("O10").Value = 0.189056558536119
$ws.Range("P10").Value = 0.189056558536119
$ws.Range("S10").Value = 0.08404839457849415
$ws.Range("T10").Value = 0.08404839457849415
